$wb = $excel.ActiveWorkbook

# --- SNVs sheet: rows 2-4 in column S (Amplicon) ---
$wsSNVs = $wb.Worksheets.Item("SNVs")
$wsSNVs.Range("S2").Value = "GRCh37.p13_chr9_88504079"
$wsSNVs.Range("S3").Value = "GRCh37.p13_chr10_89653730"
$wsSNVs.Range("S4").Value = "GRCh37.p13_chr10_89653730"

# --- Indels sheet: rows 2-41 in column S (Amplicon) ---
$wsIndels = $wb.Worksheets.Item("Indels")
for ($row = 2; $row -le 41; $row++) {
    $wsIndels.Range("S$row").Value = "GRCh37.p13_chr10_89653730"
}
